$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CustodyStatusChangeReport")

$ws.Rows.Item(52).Insert()

$ws.Range("A52").Value = ""
$ws.Range("B52").Value = "Booking Subject Location Status Description"
$ws.Range("C52").Value = "Current location status of the booking subject at the time of report"
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = "/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/j:Booking/j:BookingSubject/cscr-ext:SubjectLocationStatus/nc:StatusDescriptionText"
